{"js": "// Update the worksheet title (date) and every arithmetic-problem cell in\n// the 20x5 table with the new values, while preserving each run's\n// existing formatting (rFonts/sz) and paragraph properties (jc). We\n// replace text through an existing Range (cell.body.getRange() /\n// paragraph.getRange()) instead of clearing+inserting, so the run\n// properties and paragraph properties already on the document are kept\n// untouched and only the `w:t` text content changes, exactly like the\n// target diff.\n\nconst DATA = {\"titleOld\": \"2024-01-28 Sunday\", \"titleNew\": \"2024-01-29 Monday\", \"oldValues\": [[\"33-12=21\", \"15+28=43\", \"0+52=52\", \"79-31=48\", \"68+14=82\"], [\"24+28=52\", \"96-56=40\", \"68-23=45\", \"20+45=65\", \"66-9=57\"], [\"88-4=84\", \"77-72=5\", \"32+30=62\", \"97-31=66\", \"27-3=24\"], [\"56+13=69\", \"69-25=44\", \"69-39=30\", \"25+9=34\", \"23-8=15\"], [\"50+33=83\", \"3+91=94\", \"20-13=7\", \"59-32=27\", \"99-21=78\"], [\"39+5=44\", \"83-75=8\", \"22+3=25\", \"8+68=76\", \"45-16=29\"], [\"88-23=65\", \"35+58=93\", \"65-25=40\", \"54+0=54\", \"84-69=15\"], [\"27+69=96\", \"51-39=12\", \"95-57=38\", \"91-41=50\", \"48+22=70\"], [\"81-62=19\", \"42-8=34\", \"23+63=86\", \"45+50=95\", \"51+11=62\"], [\"73-24=49\", \"44+28=72\", \"55-11=44\", \"1+6=7\", \"46+18=64\"], [\"68-3=65\", \"36+10=46\", \"17+57=74\", \"26+48=74\", \"88-28=60\"], [\"36+23=59\", \"22+3=25\", \"38+9=47\", \"5+69=74\", \"72-55=17\"], [\"6+82=88\", \"56-35=21\", \"67-39=28\", \"36-18=18\", \"58+19=77\"], [\"21+75=96\", \"33-31=2\", \"40+38=78\", \"18-14=4\", \"82-21=61\"], [\"89+9=98\", \"28+35=63\", \"87-20=67\", \"96-45=51\", \"75-12=63\"], [\"51+1=52\", \"99-95=4\", \"45+19=64\", \"52+0=52\", \"28+49=77\"], [\"25-3=22\", \"17+40=57\", \"64+21=85\", \"58-18=40\", \"22+6=28\"], [\"69-47=22\", \"65+21=86\", \"81-15=66\", \"0+82=82\", \"7+64=71\"], [\"32-20=12\", \"45+7=52\", \"43+32=75\", \"51+12=63\", \"9-7=2\"], [\"3+59=62\", \"59+19=78\", \"22-14=8\", \"83-45=38\", \"71+21=92\"]], \"newValues\": [[\"25-17=8\", \"84-6=78\", \"92-50=42\", \"99-54=45\", \"9+38=47\"], [\"64-11=53\", \"34+17=51\", \"94-18=76\", \"15-2=13\", \"78-52=26\"], [\"4+82=86\", \"59-39=20\", \"23-4=19\", \"45+49=94\", \"1+25=26\"], [\"98-56=42\", \"62-33=29\", \"32+53=85\", \"25+3=28\", \"94-90=4\"], [\"78+0=78\", \"48+20=68\", \"91-22=69\", \"80-41=39\", \"25-11=14\"], [\"64-20=44\", \"80-77=3\", \"50+21=71\", \"12+21=33\", \"57-42=15\"], [\"84-44=40\", \"60-24=36\", \"22+55=77\", \"48-2=46\", \"30+53=83\"], [\"29+20=49\", \"62-3=59\", \"78-55=23\", \"46-8=38\", \"96-39=57\"], [\"91-85=6\", \"3-3=0\", \"76+17=93\", \"93+2=95\", \"99-98=1\"], [\"98-94=4\", \"54-14=40\", \"33+20=53\", \"28+19=47\", \"77-11=66\"], [\"41+25=66\", \"73-16=57\", \"57+2=59\", \"30+30=60\", \"35+20=55\"], [\"31+43=74\", \"39+59=98\", \"46-30=16\", \"89+3=92\", \"48+6=54\"], [\"74-68=6\", \"84-40=44\", \"84-26=58\", \"90-26=64\", \"18+31=49\"], [\"51+6=57\", \"96-30=66\", \"97-84=13\", \"36+15=51\", \"49+9=58\"], [\"1+47=48\", \"77+0=77\", \"30+31=61\", \"43+2=45\", \"2+42=44\"], [\"77+7=84\", \"15+11=26\", \"17+29=46\", \"53-40=13\", \"19+66=85\"], [\"48-2=46\", \"54+30=84\", \"2+48=50\", \"80-66=14\", \"18-15=3\"], [\"45+40=85\", \"63-31=32\", \"19+3=22\", \"47+23=70\", \"10+76=86\"], [\"18+45=63\", \"72+15=87\", \"51-2=49\", \"66-39=27\", \"45-23=22\"], [\"70-3=67\", \"18+74=92\", \"70+22=92\", \"78+21=99\", \"62-19=43\"]]};\n\nconst body = context.document.body;\n\n// --- Title paragraph (first paragraph in the body holds the date) ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.load(\"text\");\nawait context.sync();\n\nif (titlePara.text.trim() !== DATA.titleOld) {\n  throw new Error(\n    \"Unexpected title text: \" + JSON.stringify(titlePara.text) +\n    \" (expected \" + JSON.stringify(DATA.titleOld) + \")\"\n  );\n}\ntitlePara.getRange().insertText(DATA.titleNew, Word.InsertLocation.replace);\n\n// --- Table cells (20 rows x 5 columns of arithmetic problems) ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst oldValues = DATA.oldValues;\nconst newValues = DATA.newValues;\n\nfor (let r = 0; r < newValues.length; r++) {\n  const oldRow = oldValues[r];\n  const newRow = newValues[r];\n  for (let c = 0; c < newRow.length; c++) {\n    const cell = table.getCell(r, c);\n    const currentValue = table.values[r][c];\n    if (currentValue !== oldRow[c]) {\n      throw new Error(\n        \"Unexpected cell (\" + r + \",\" + c + \") text: \" +\n        JSON.stringify(currentValue) + \" (expected \" + JSON.stringify(oldRow[c]) + \")\"\n      );\n    }\n    cell.body.getRange().insertText(newRow[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet title (date) and every arithmetic-problem cell in\n# the 20x5 table with the new values. Setting Range.Text in place (rather\n# than deleting + re-inserting a run) keeps each run's existing formatting\n# (rFonts/sz) and paragraph properties (jc) untouched, matching the target\n# diff which only changes the `w:t` text content.\n\n$d = $word.ActiveDocument\n\n# --- Title paragraph (first paragraph in the body holds the date) ---\n$titlePara = $d.Paragraphs.Item(1)\n$expectedTitle = \"2024-01-28 Sunday\"\n$actualTitle = $titlePara.Range.Text.TrimEnd([char]13, [char]7)\nif ($actualTitle -ne $expectedTitle) {\n    throw \"Unexpected title text: [$actualTitle] (expected [$expectedTitle])\"\n}\n$titlePara.Range.Text = \"2024-01-29 Monday\"\n\n# --- Table cells (20 rows x 5 columns of arithmetic problems) ---\n$tbl = $d.Tables.Item(1)\n\n$oldValues = @(\n    ,@(\"33-12=21\", \"15+28=43\", \"0+52=52\", \"79-31=48\", \"68+14=82\")\n    ,@(\"24+28=52\", \"96-56=40\", \"68-23=45\", \"20+45=65\", \"66-9=57\")\n    ,@(\"88-4=84\", \"77-72=5\", \"32+30=62\", \"97-31=66\", \"27-3=24\")\n    ,@(\"56+13=69\", \"69-25=44\", \"69-39=30\", \"25+9=34\", \"23-8=15\")\n    ,@(\"50+33=83\", \"3+91=94\", \"20-13=7\", \"59-32=27\", \"99-21=78\")\n    ,@(\"39+5=44\", \"83-75=8\", \"22+3=25\", \"8+68=76\", \"45-16=29\")\n    ,@(\"88-23=65\", \"35+58=93\", \"65-25=40\", \"54+0=54\", \"84-69=15\")\n    ,@(\"27+69=96\", \"51-39=12\", \"95-57=38\", \"91-41=50\", \"48+22=70\")\n    ,@(\"81-62=19\", \"42-8=34\", \"23+63=86\", \"45+50=95\", \"51+11=62\")\n    ,@(\"73-24=49\", \"44+28=72\", \"55-11=44\", \"1+6=7\", \"46+18=64\")\n    ,@(\"68-3=65\", \"36+10=46\", \"17+57=74\", \"26+48=74\", \"88-28=60\")\n    ,@(\"36+23=59\", \"22+3=25\", \"38+9=47\", \"5+69=74\", \"72-55=17\")\n    ,@(\"6+82=88\", \"56-35=21\", \"67-39=28\", \"36-18=18\", \"58+19=77\")\n    ,@(\"21+75=96\", \"33-31=2\", \"40+38=78\", \"18-14=4\", \"82-21=61\")\n    ,@(\"89+9=98\", \"28+35=63\", \"87-20=67\", \"96-45=51\", \"75-12=63\")\n    ,@(\"51+1=52\", \"99-95=4\", \"45+19=64\", \"52+0=52\", \"28+49=77\")\n    ,@(\"25-3=22\", \"17+40=57\", \"64+21=85\", \"58-18=40\", \"22+6=28\")\n    ,@(\"69-47=22\", \"65+21=86\", \"81-15=66\", \"0+82=82\", \"7+64=71\")\n    ,@(\"32-20=12\", \"45+7=52\", \"43+32=75\", \"51+12=63\", \"9-7=2\")\n    ,@(\"3+59=62\", \"59+19=78\", \"22-14=8\", \"83-45=38\", \"71+21=92\")\n)\n\n$newValues = @(\n    ,@(\"25-17=8\", \"84-6=78\", \"92-50=42\", \"99-54=45\", \"9+38=47\")\n    ,@(\"64-11=53\", \"34+17=51\", \"94-18=76\", \"15-2=13\", \"78-52=26\")\n    ,@(\"4+82=86\", \"59-39=20\", \"23-4=19\", \"45+49=94\", \"1+25=26\")\n    ,@(\"98-56=42\", \"62-33=29\", \"32+53=85\", \"25+3=28\", \"94-90=4\")\n    ,@(\"78+0=78\", \"48+20=68\", \"91-22=69\", \"80-41=39\", \"25-11=14\")\n    ,@(\"64-20=44\", \"80-77=3\", \"50+21=71\", \"12+21=33\", \"57-42=15\")\n    ,@(\"84-44=40\", \"60-24=36\", \"22+55=77\", \"48-2=46\", \"30+53=83\")\n    ,@(\"29+20=49\", \"62-3=59\", \"78-55=23\", \"46-8=38\", \"96-39=57\")\n    ,@(\"91-85=6\", \"3-3=0\", \"76+17=93\", \"93+2=95\", \"99-98=1\")\n    ,@(\"98-94=4\", \"54-14=40\", \"33+20=53\", \"28+19=47\", \"77-11=66\")\n    ,@(\"41+25=66\", \"73-16=57\", \"57+2=59\", \"30+30=60\", \"35+20=55\")\n    ,@(\"31+43=74\", \"39+59=98\", \"46-30=16\", \"89+3=92\", \"48+6=54\")\n    ,@(\"74-68=6\", \"84-40=44\", \"84-26=58\", \"90-26=64\", \"18+31=49\")\n    ,@(\"51+6=57\", \"96-30=66\", \"97-84=13\", \"36+15=51\", \"49+9=58\")\n    ,@(\"1+47=48\", \"77+0=77\", \"30+31=61\", \"43+2=45\", \"2+42=44\")\n    ,@(\"77+7=84\", \"15+11=26\", \"17+29=46\", \"53-40=13\", \"19+66=85\")\n    ,@(\"48-2=46\", \"54+30=84\", \"2+48=50\", \"80-66=14\", \"18-15=3\")\n    ,@(\"45+40=85\", \"63-31=32\", \"19+3=22\", \"47+23=70\", \"10+76=86\")\n    ,@(\"18+45=63\", \"72+15=87\", \"51-2=49\", \"66-39=27\", \"45-23=22\")\n    ,@(\"70-3=67\", \"18+74=92\", \"70+22=92\", \"78+21=99\", \"62-19=43\")\n)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n    $oldRow = $oldValues[$r]\n    $newRow = $newValues[$r]\n    for ($c = 0; $c -lt $newRow.Count; $c++) {\n        $cell = $tbl.Cell($r + 1, $c + 1)\n        $actual = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($actual -ne $oldRow[$c]) {\n            throw \"Unexpected cell ($r,$c) text: [$actual] (expected [$($oldRow[$c])])\"\n        }\n        $cell.Range.Text = $newRow[$c]\n    }\n}\n\n"}
